$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("F2").Value = 3.75
$ws.Range("G2").Value = 3.85
$ws.Range("H2").Value = 2.06
$ws.Range("I2").Value = 2.08
$ws.Range("J2").Value = 3.9
$ws.Range("K2").Value = 4
$ws.Range("O2").Value = 1.32
$ws.Range("Q2").Value = 1.93
$ws.Range("S2").Value = 3.4
$ws.Range("V2").Value = 1.92
$ws.Range("W2").Value = 1.35
$ws.Range("AB2").Value = 15
$ws.Range("AF2").Value = 29
$ws.Range("AG2").Value = 16
$ws.Range("AL2").Value = 55

# Row 3
$ws.Range("F3").Value = 2.52
$ws.Range("G3").Value = 2.7
$ws.Range("H3").Value = 2.64
$ws.Range("I3").Value = 2.82
$ws.Range("N3").Value = 5.7
$ws.Range("P3").Value = 2.58
$ws.Range("T3").Value = 1.54
$ws.Range("U3").Value = 2.7
$ws.Range("V3").Value = 1.55
$ws.Range("W3").Value = 1.59
$ws.Range("AG3").Value = 12.5
$ws.Range("AN3").Value = 15.5
$ws.Range("AO3").Value = 15

# Row 4
$ws.Range("Q4").Value = 1.83

# Row 5
$ws.Range("F5").Value = 1.04
$ws.Range("H5").Value = 1.13
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 2.52
$ws.Range("K5").Value = 950

# Row 6
$ws.Range("P6").Value = 1.86
$ws.Range("Q6").Value = 1.84

# Row 7
$ws.Range("F7").Value = 9.199999999999999
$ws.Range("G7").Value = 14
$ws.Range("H7").Value = 1.25
$ws.Range("J7").Value = 6.6
$ws.Range("K7").Value = 9.6

# Row 9
$ws.Range("F9").Value = 3.1
$ws.Range("G9").Value = 3.4
$ws.Range("H9").Value = 2.7
$ws.Range("I9").Value = 2.96
$ws.Range("J9").Value = 2.86
$ws.Range("K9").Value = 3.15

# Row 12
$ws.Range("O12").Value = 1.6
$ws.Range("T12").Value = 2.3
$ws.Range("AE12").Value = 1000
$ws.Range("AL12").Value = 65

# Row 13
$ws.Range("P13").Value = 1.56
$ws.Range("AD13").Value = 21
$ws.Range("AO13").Value = 170
